$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "File Name" values (column H) for the existing rows ---
# (entered top-to-bottom, which is also the order the strings were first introduced)
$ws.Range("H5").Value = "protectYourself.html"
$ws.Range("H6").Value = "amazonHP.html"
$ws.Range("H7").Value = "irs.html"
$ws.Range("H8").Value = "IRS_YourAccount.html"
$ws.Range("H9").Value = "appleID.html"
$ws.Range("H10").Value = "amazonProductPage.html"
$ws.Range("H11").Value = "mcafee.html"
$ws.Range("H12").Value = "getProtected.html"
$ws.Range("H14").Value = "redcross_covidRelief.html"
$ws.Range("H15").Value = "ssa.html"
$ws.Range("H16").Value = "mySocialSecurity.html"
$ws.Range("H17").Value = "amazon_maskDelivery.html"
$ws.Range("H18").Value = "ssaFb.html"
$ws.Range("H19").Value = "ssa_optOut.html"
$ws.Range("H20").Value = "ssa_replacementCard.html"
$ws.Range("H21").Value = "walmart.html"

# --- Bug fix: row 18's Imposter Type should be Biz+Govt, not Govt ---
$ws.Range("F18").Value = "Biz+Govt"

# --- New rows for the Letter-type content ---
$ws.Range("B22").Value = "Letter"
$ws.Range("C22").Value = "Test"
$ws.Range("D22").Value = "Real"
$ws.Range("F22").Value = "Govt"
$ws.Range("H22").Value = "medicareReview"

$ws.Range("B23").Value = "Letter"
$ws.Range("C23").Value = "Test"
$ws.Range("D23").Value = "Scam"
$ws.Range("F23").Value = "Govt"
$ws.Range("H23").Value = "benefitsSuspension"

# --- New "Tooltips?" column (J) ---
$ws.Range("J3").Value = "Tooltips?"
$ws.Range("G3").Font.Bold = $true

$ws.Range("J5").Value = "Yes"
$ws.Range("J6").Value = "yes"
$ws.Range("J7").Value = "Yes"
$ws.Range("J8").Value = "Yes"
$ws.Range("J9").Value = "yes"
$ws.Range("J10").Value = "Yes"
$ws.Range("J11").Value = "Yes"
$ws.Range("J12").Value = "yes"

# --- Column width tweaks (narrower B/C/E columns) ---
$ws.Columns("B").ColumnWidth = 6.6
$ws.Columns("C").ColumnWidth = 10.26
$ws.Columns("E").ColumnWidth = 19.26

# --- Selection / scroll position left by the author when they saved ---
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H17").Select() | Out-Null
